$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Sending cluster): ECs -> MuSCs for data rows 2-4 ---
$ws.Range("A2").Value = "MuSCs"
$ws.Range("A3").Value = "MuSCs"
$ws.Range("A4").Value = "MuSCs"

# --- Row 2 (Target cluster = ECs) ---
$ws.Range("G2").Value = 0.01490866666666667
$ws.Range("H2").Value = 0.044726
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 26.532132
$ws.Range("N2").Value = 79.596396
$ws.Range("O2").Value = 0.3960736634233649
$ws.Range("P2").Value = 0.3960736634233648
$ws.Range("Q2").Value = 0.3955587119440001
$ws.Range("R2").Value = 3.560028407496
$ws.Range("S2").Value = 0.3960736634233649
$ws.Range("T2").Value = 0.3960736634233648

# --- Row 3 (Target cluster = FAPs) ---
$ws.Range("G3").Value = 0.01490866666666667
$ws.Range("H3").Value = 0.044726
$ws.Range("O3").Value = 0.2505213219764053
$ws.Range("P3").Value = 0.2505213219764053
$ws.Range("Q3").Value = 0.2501956090162222
$ws.Range("R3").Value = 2.251760481146
$ws.Range("S3").Value = 0.2505213219764053
$ws.Range("T3").Value = 0.2505213219764053

# --- Row 4 (Target cluster = MuSCs) ---
$ws.Range("G4").Value = 0.01490866666666667
$ws.Range("H4").Value = 0.044726
$ws.Range("M4").Value = 23.67385
$ws.Range("N4").Value = 71.02154999999999
$ws.Range("O4").Value = 0.3534050146002298
$ws.Range("P4").Value = 0.3534050146002298
$ws.Range("Q4").Value = 0.3529455383666666
$ws.Range("R4").Value = 3.1765098453
$ws.Range("S4").Value = 0.3534050146002298
$ws.Range("T4").Value = 0.3534050146002298
